# Generate Report for Handoff
# Inserts a new tracked file "5fb35e66-d076-4211-8a4f-a1c4607d1cc9.md" into the
# localization-status report, right before the existing
# "da9e7863-03e8-4548-a100-4732be3f6675.md" entry (which was previously the
# last row) on all three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$mainHash  = "9ff4fb3b38248654e208e96ea57850f1f78f4178"
$zhHash    = "93d01c20513f765c884212f724569292b88ab7bb"
$deHash    = "93d01c20513f765c884212f724569292b88ab7bb"
$newFile   = "5fb35e66-d076-4211-8a4f-a1c4607d1cc9.md"
$newDateMain = "2016-09-08 04:58:36"
$newDateZh   = "2016-09-08 04:58:31"
$newDateDe   = "2016-09-08 04:58:36"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")

$wsO.Rows.Item(8).Insert()

$wsO.Range("A8").Value = $newFile
$wsO.Range("B8").Value = "e2e\" + $newFile
$wsO.Range("C8").Value = ".md"
$wsO.Range("D8").Value = ""
$wsO.Range("E8").Value = "Ready for handoff"
$wsO.Range("F8").Value = "Ready for handoff"
$wsO.Range("G8").Value = $newDateMain

$loO = $wsO.ListObjects.Item("Overview")
$loO.Resize($wsO.Range("A1:G9"))

$wsO.Range("A1").Hyperlinks.Delete()
$null = $wsO.Hyperlinks.Add($wsO.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe1cd128624bfdeab48986cfcc4c49ee29e92d24/e2e/86c14bdc-2d76-44ad-bb1b-cbe6d32d2268.md", "", "", "e2e\86c14bdc-2d76-44ad-bb1b-cbe6d32d2268.md")
$null = $wsO.Hyperlinks.Add($wsO.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e963374495c2bce247d180056ededf68b165dcdb/e2e/2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md", "", "", "e2e\2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md")
$null = $wsO.Hyperlinks.Add($wsO.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3696c9ac59a4fade26a5a2f728fdac23ee7f628/e2e/34bfdbb5-28e8-43bc-b318-a409e2af5021.md", "", "", "e2e\34bfdbb5-28e8-43bc-b318-a409e2af5021.md")
$null = $wsO.Hyperlinks.Add($wsO.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e963374495c2bce247d180056ededf68b165dcdb/e2e/61a15d7d-b6d5-4da7-b456-f9204bdc3269.md", "", "", "e2e\61a15d7d-b6d5-4da7-b456-f9204bdc3269.md")
$null = $wsO.Hyperlinks.Add($wsO.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fa444c4e6e4c9526e4ae0053a8525b895f1e0a29/e2e/e5a52eae-b88d-463c-842a-daaf45f1639a.md", "", "", "e2e\e5a52eae-b88d-463c-842a-daaf45f1639a.md")
$null = $wsO.Hyperlinks.Add($wsO.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad83d812831d97b9452e9c2d81804c936654ff90/e2e/059e2db7-ca47-4434-a3e6-e74353cdbddd.md", "", "", "e2e\059e2db7-ca47-4434-a3e6-e74353cdbddd.md")
$null = $wsO.Hyperlinks.Add($wsO.Range("B8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$mainHash/e2e/$newFile", "", "", "e2e\" + $newFile)
$null = $wsO.Hyperlinks.Add($wsO.Range("B9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb52c28d5f7d7afe77c667df0efda4c9903844fc/e2e/da9e7863-03e8-4548-a100-4732be3f6675.md", "", "", "e2e\da9e7863-03e8-4548-a100-4732be3f6675.md")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZ = $wb.Worksheets.Item("zh-cn")

$wsZ.Rows.Item(8).Insert()

$wsZ.Range("A8").Value = $newFile
$wsZ.Range("B8").Value = ".md"
$wsZ.Range("C8").Value = "Ready for handoff"
$wsZ.Range("D8").Value = "e2e"
$wsZ.Range("E8").Value = "ht"
$wsZ.Range("F8").Value = "False"
$wsZ.Range("G8").Value = "$newFile.$zhHash.zh-cn.xlf"
$wsZ.Range("H8").Value = $newDateZh
$wsZ.Range("I8").Value = ""
$wsZ.Range("J8").Value = ""
$wsZ.Range("K8").Value = "0001-01-01 00:00:00"
$wsZ.Range("L8").Value = ""
$wsZ.Range("M8").Value = "True"
$wsZ.Range("N8").Value = ""
$wsZ.Range("O8").Value = "False"
$wsZ.Range("P8").Value = ""

$loZ = $wsZ.ListObjects.Item("zh-cn")
$loZ.Resize($wsZ.Range("A1:P9"))

$wsZ.Range("A1").Hyperlinks.Delete()
$null = $wsZ.Hyperlinks.Add($wsZ.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe1cd128624bfdeab48986cfcc4c49ee29e92d24/e2e/86c14bdc-2d76-44ad-bb1b-cbe6d32d2268.md", "", "", "86c14bdc-2d76-44ad-bb1b-cbe6d32d2268.md")
$null = $wsZ.Hyperlinks.Add($wsZ.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d40a11274015af4bfedc649f01fadc19ef19c48a/e2e/86c14bdc-2d76-44ad-bb1b-cbe6d32d2268.md", "", "", "86c14bdc-2d76-44ad-bb1b-cbe6d32d2268.md")
$null = $wsZ.Hyperlinks.Add($wsZ.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e963374495c2bce247d180056ededf68b165dcdb/e2e/2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md", "", "", "2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md")
$null = $wsZ.Hyperlinks.Add($wsZ.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3696c9ac59a4fade26a5a2f728fdac23ee7f628/e2e/34bfdbb5-28e8-43bc-b318-a409e2af5021.md", "", "", "34bfdbb5-28e8-43bc-b318-a409e2af5021.md")
$null = $wsZ.Hyperlinks.Add($wsZ.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e963374495c2bce247d180056ededf68b165dcdb/e2e/61a15d7d-b6d5-4da7-b456-f9204bdc3269.md", "", "", "61a15d7d-b6d5-4da7-b456-f9204bdc3269.md")
$null = $wsZ.Hyperlinks.Add($wsZ.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fa444c4e6e4c9526e4ae0053a8525b895f1e0a29/e2e/e5a52eae-b88d-463c-842a-daaf45f1639a.md", "", "", "e5a52eae-b88d-463c-842a-daaf45f1639a.md")
$null = $wsZ.Hyperlinks.Add($wsZ.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/bbe484687a5672796c28c94868bcb1a5dac80c78/e2e/e5a52eae-b88d-463c-842a-daaf45f1639a.md", "", "", "e5a52eae-b88d-463c-842a-daaf45f1639a.md")
$null = $wsZ.Hyperlinks.Add($wsZ.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad83d812831d97b9452e9c2d81804c936654ff90/e2e/059e2db7-ca47-4434-a3e6-e74353cdbddd.md", "", "", "059e2db7-ca47-4434-a3e6-e74353cdbddd.md")
$null = $wsZ.Hyperlinks.Add($wsZ.Range("A8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$mainHash/e2e/$newFile", "", "", $newFile)
$null = $wsZ.Hyperlinks.Add($wsZ.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb52c28d5f7d7afe77c667df0efda4c9903844fc/e2e/da9e7863-03e8-4548-a100-4732be3f6675.md", "", "", "da9e7863-03e8-4548-a100-4732be3f6675.md")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsD = $wb.Worksheets.Item("de-de")

$wsD.Rows.Item(8).Insert()

$wsD.Range("A8").Value = $newFile
$wsD.Range("B8").Value = ".md"
$wsD.Range("C8").Value = "Ready for handoff"
$wsD.Range("D8").Value = "e2e"
$wsD.Range("E8").Value = "ht"
$wsD.Range("F8").Value = "False"
$wsD.Range("G8").Value = "$newFile.$deHash.de-de.xlf"
$wsD.Range("H8").Value = $newDateDe
$wsD.Range("I8").Value = ""
$wsD.Range("J8").Value = ""
$wsD.Range("K8").Value = "0001-01-01 00:00:00"
$wsD.Range("L8").Value = ""
$wsD.Range("M8").Value = "True"
$wsD.Range("N8").Value = ""
$wsD.Range("O8").Value = "False"
$wsD.Range("P8").Value = ""

$loD = $wsD.ListObjects.Item("de-de")
$loD.Resize($wsD.Range("A1:P9"))

$wsD.Range("A1").Hyperlinks.Delete()
$null = $wsD.Hyperlinks.Add($wsD.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe1cd128624bfdeab48986cfcc4c49ee29e92d24/e2e/86c14bdc-2d76-44ad-bb1b-cbe6d32d2268.md", "", "", "86c14bdc-2d76-44ad-bb1b-cbe6d32d2268.md")
$null = $wsD.Hyperlinks.Add($wsD.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ad8d7faae85f56ae966965a9a1a121d72e4ff6ac/e2e/86c14bdc-2d76-44ad-bb1b-cbe6d32d2268.md", "", "", "86c14bdc-2d76-44ad-bb1b-cbe6d32d2268.md")
$null = $wsD.Hyperlinks.Add($wsD.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e963374495c2bce247d180056ededf68b165dcdb/e2e/2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md", "", "", "2bdc51fe-c985-4dcb-995d-f1ac8ee5abb8.md")
$null = $wsD.Hyperlinks.Add($wsD.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3696c9ac59a4fade26a5a2f728fdac23ee7f628/e2e/34bfdbb5-28e8-43bc-b318-a409e2af5021.md", "", "", "34bfdbb5-28e8-43bc-b318-a409e2af5021.md")
$null = $wsD.Hyperlinks.Add($wsD.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e963374495c2bce247d180056ededf68b165dcdb/e2e/61a15d7d-b6d5-4da7-b456-f9204bdc3269.md", "", "", "61a15d7d-b6d5-4da7-b456-f9204bdc3269.md")
$null = $wsD.Hyperlinks.Add($wsD.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fa444c4e6e4c9526e4ae0053a8525b895f1e0a29/e2e/e5a52eae-b88d-463c-842a-daaf45f1639a.md", "", "", "e5a52eae-b88d-463c-842a-daaf45f1639a.md")
$null = $wsD.Hyperlinks.Add($wsD.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4d2c9a8534dc4a6aab9285f5b0285e7f31e2a651/e2e/e5a52eae-b88d-463c-842a-daaf45f1639a.md", "", "", "e5a52eae-b88d-463c-842a-daaf45f1639a.md")
$null = $wsD.Hyperlinks.Add($wsD.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad83d812831d97b9452e9c2d81804c936654ff90/e2e/059e2db7-ca47-4434-a3e6-e74353cdbddd.md", "", "", "059e2db7-ca47-4434-a3e6-e74353cdbddd.md")
$null = $wsD.Hyperlinks.Add($wsD.Range("A8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$mainHash/e2e/$newFile", "", "", $newFile)
$null = $wsD.Hyperlinks.Add($wsD.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb52c28d5f7d7afe77c667df0efda4c9903844fc/e2e/da9e7863-03e8-4548-a100-4732be3f6675.md", "", "", "da9e7863-03e8-4548-a100-4732be3f6675.md")
